$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "72.476.83"
$ws.Range("E2").Value = "  +5.70%  "

# Row 3
$ws.Range("D3").Value = "4.069.26"
$ws.Range("E3").Value = "  +6.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.41"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.80"
$ws.Range("E6").Value = "  +3.60%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.727"
$ws.Range("E7").Value = "  +20.69%  "

# Row 8
$ws.Range("D8").Value = "4.062.27"
$ws.Range("E8").Value = "  +6.05%  "

# Row 9
$ws.Range("E9").Value = "  +0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.784"
$ws.Range("E10").Value = "  +10.75%  "

# Row 11
$ws.Range("E11").Value = "  +6.04%  "

# Row 12
$ws.Range("E12").Value = "  +2.43%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.94"
$ws.Range("E13").Value = "  +17.98%  "

# Row 14
$ws.Range("E14").Value = "  +10.23%  "

# Row 15
$ws.Range("D15").Value = "4.707.95"
$ws.Range("E15").Value = "  +6.16%  "

# Row 16
$ws.Range("D16").Value = "4.081.44"
$ws.Range("E16").Value = "  +6.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.53"
$ws.Range("E17").Value = "  +5.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "21.44"
$ws.Range("E18").Value = "  +2.68%  "

# Row 19
$ws.Range("E19").Value = "  +2.47%  "

# Row 20
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("D21").Value = "72.497.51"
$ws.Range("E21").Value = "  +5.85%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.68"
$ws.Range("E22").Value = "  +7.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "104.34"
$ws.Range("E23").Value = "  +21.01%  "

# Row 24
$ws.Range("E24").Value = "  +6.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.09"
$ws.Range("E25").Value = "  +8.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  +2.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("E27").Value = "  +0.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.14"
$ws.Range("E28").Value = "  +6.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.19"
$ws.Range("E29").Value = "  +6.57%  "

# Row 30
$ws.Range("E30").Value = "  +3.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.29"
$ws.Range("E31").Value = "  +16.81%  "

# Row 32
$ws.Range("E32").Value = "  +5.43%  "

# Row 33
$ws.Range("E33").Value = "  +5.23%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "682.76"
$ws.Range("E34").Value = "  +0.21%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.95"
$ws.Range("E35").Value = "  +1.59%  "

# Row 36
$ws.Range("E36").Value = "  +13.53%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.28"
$ws.Range("E37").Value = "  +7.35%  "

# Row 38
$ws.Range("E38").Value = "  +2.44%  "

# Row 39
$ws.Range("E39").Value = "  -0.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.153"
$ws.Range("E40").Value = "  +5.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.48"
$ws.Range("E41").Value = "  +9.39%  "

# Row 42
$ws.Range("E42").Value = "  +0.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0501"
$ws.Range("E43").Value = "  +5.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").Value = "  +3.03%  "

# Row 46
$ws.Range("E46").Value = "  +14.77%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.88"
$ws.Range("E47").Value = "  +17.94%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.44"
$ws.Range("E49").Value = "  +1.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.09"
$ws.Range("E50").Value = "  +5.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000276"
$ws.Range("E51").Value = "  +2.21%  "

